# Corrections following third round of review:
# Remove the "subgenus" field/column from the Materials sheet.
# This deletes the column whose header is "subgenus" (and whose row-2
# value is the template placeholder "${subgenus}"), shifting every
# column to its right one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Locate the "subgenus" header column dynamically (falls back to the
# known column letter "AS" if, for some reason, it cannot be found).
$headerCell = $ws.Rows(1).Find("subgenus")
if ($headerCell -ne $null) {
    $targetColumn = $headerCell.Column
} else {
    $targetColumn = 45
}

$ws.Columns($targetColumn).Delete()
